$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 38

$ws.Cells.Item($row, 1).Value = "71277620"
$ws.Cells.Item($row, 2).Value = ""
$ws.Cells.Item($row, 3).Value = "Cash"
$ws.Cells.Item($row, 4).Value = "2025-08-18T17:28:47"
$ws.Cells.Item($row, 5).Value = 76
$ws.Cells.Item($row, 6).Value = ""
$ws.Cells.Item($row, 7).Value = 76
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
